{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the two target paragraphs by their (pre-edit) concatenated text.\nlet eulersPara = null;\nlet skyboxPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"Used eulers angles\") {\n    eulersPara = paragraphs.items[i];\n  } else if (t === \"Started with skybox class, then re-written it with the cubemap class \") {\n    skyboxPara = paragraphs.items[i];\n  }\n}\n\nif (!eulersPara || !skyboxPara) {\n  throw new Error(\"Could not locate expected paragraphs to edit.\");\n}\n\n// 1) Collapse \"Used \"/\"eulers\"/\" angles\" (split across runs because of the\n//    spell-check proofErr wrapper around \"eulers\") into a single clean run.\neulersPara.insertText(\"Used eulers angles\", \"Replace\");\n\n// 2) Same cleanup for the skybox paragraph: merge its three runs (the\n//    proofErr-wrapped \"cubemap\" run included) back into one run.\nskyboxPara.insertText(\n  \"Started with skybox class, then re-written it with the cubemap class \",\n  \"Replace\"\n);\nawait context.sync();\n\n// 3) Add the new \"Airplane texture...\" paragraph right after it. Insert it\n//    relative to the next (blank, non-list) paragraph so the new paragraph\n//    picks up plain body formatting instead of inheriting the bullet/list\n//    formatting of the skybox paragraph.\nconst nextPara = skyboxPara.getNext();\nnextPara.insertParagraph(\"Airplane texture didn\\u2019t seem to work properly..\", \"Before\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the two target paragraphs by their (pre-edit) concatenated text so\n# this does not depend on hard-coded paragraph indices.\n$eulersIndex = 0\n$skyboxIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -eq \"Used eulers angles`r\") {\n        $eulersIndex = $i\n    } elseif ($t -eq \"Started with skybox class, then re-written it with the cubemap class `r\") {\n        $skyboxIndex = $i\n    }\n}\n\nif ($eulersIndex -eq 0 -or $skyboxIndex -eq 0) {\n    throw \"Could not locate expected paragraphs to edit.\"\n}\n\n# 1) Collapse \"Used \"/\"eulers\"/\" angles\" (split across runs because of the\n#    spell-check proofErr wrapper around \"eulers\") into a single clean run.\n$p1 = $d.Paragraphs.Item($eulersIndex)\n$r1 = $p1.Range\n$r1.MoveEnd(1, -1) | Out-Null          # exclude the paragraph mark\n$lang1 = $r1.LanguageID\n$r1.Delete()\n$r1.InsertAfter(\"Used eulers angles\")\n$r1.LanguageID = $lang1\n\n# 2) Same cleanup for the skybox paragraph: merge its three runs (the\n#    proofErr-wrapped \"cubemap\" run included) back into one run.\n$p2 = $d.Paragraphs.Item($skyboxIndex)\n$r2 = $p2.Range\n$r2.MoveEnd(1, -1) | Out-Null\n$lang2 = $r2.LanguageID\n$r2.Delete()\n$r2.InsertAfter(\"Started with skybox class, then re-written it with the cubemap class \")\n$r2.LanguageID = $lang2\n\n# 3) Add the new \"Airplane texture...\" paragraph right after it. Insert the\n#    paragraph break before the next (blank, non-list) paragraph so the new\n#    paragraph picks up plain body formatting instead of inheriting the\n#    bullet/list formatting of the skybox paragraph.\n$nextPara = $d.Paragraphs.Item($skyboxIndex + 1)\n$nextRange = $nextPara.Range\n$nextRange.Collapse(1) | Out-Null      # wdCollapseStart\n$nextRange.InsertParagraphBefore()\n\n$newPara = $d.Paragraphs.Item($skyboxIndex + 1)\n$newPara.Range.InsertBefore(\"Airplane texture didn\" + [char]0x2019 + \"t seem to work properly..\")\n"}
